# The edit lives on the "2020-11-21" attendance sheet: a new row 6 is
# appended below the existing data (rows 1-5), reusing the same
# bold/bordered/centered formatting as the rest of column A, and the
# sheet's used-range grows from A1:I5 to A1:I6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020-11-21")

# Clone row 5's formatting (bold font, thin border, centered alignment)
# into row 6's Sr. No cell, the same way the source row's look carries
# forward when a new attendance entry is appended.
$ws.Range("A5").Copy($ws.Range("A6"))

$ws.Range("B6").Value = "sachin"
$ws.Range("C6").Value = "301/Sanskruti-1,Andheri, Mumbai"
$ws.Range("D6").Value = "Software Engineer"
$ws.Range("E6").Value = "15:41:33"
$ws.Range("F6").Value = 93.81816225047248
$ws.Range("G6").Value = 55.61066115619608
$ws.Range("H6").Value = "NA"
$ws.Range("I6").Value = "NA"
